# Fruta / hortaliza, semanal
# Insert a new weekly record before the existing row 205 (old data shifts
# down by one row, and the former last row re-appears one row lower with a
# brand-new row appended at the very end of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 205:322 down to 206:323, opening up a blank row 205.
$ws.Rows.Item(205).Insert()

# Populate the newly opened row 205 with the new weekly observation.
$ws.Cells.Item(205, 1).Value  = 7
$ws.Cells.Item(205, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(205, 3).Value  = "Ñuble"
$ws.Cells.Item(205, 4).Value  = 44806
$ws.Cells.Item(205, 5).Value  = 16
$ws.Cells.Item(205, 6).Value  = 100114013
$ws.Cells.Item(205, 7).Value  = "Zanahoria"
$ws.Cells.Item(205, 8).Value  = "Sin especificar"
$ws.Cells.Item(205, 9).Value  = "Primera"
$ws.Cells.Item(205, 10).Value = 120
$ws.Cells.Item(205, 11).Value = 9000
$ws.Cells.Item(205, 12).Value = 9500
$ws.Cells.Item(205, 13).Value = 9250
$ws.Cells.Item(205, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(205, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(205, 16).Value = 462
$ws.Cells.Item(205, 17).Value = 20
$ws.Cells.Item(205, 18).Value = "Hortaliza"
